# Append: 2026-01-11 02:09 JST
# Update the "取得日時" (acquired datetime) column on the ランサーズ sheet
# for all existing data rows (2-10) to the new run's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-11 02:09:42"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
